$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while preserving its original "no explicit
# style" (index 0) formatting, even for values that look numeric (so Excel
# doesn't silently convert them to a stored number). We briefly force a text
# number format so the auto-type-detection leaves the string alone, then
# restore the cell to the "Normal" style (index 0), matching the workbook's
# original per-cell styling.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row -> new D/E values ($null means "leave this one unchanged")
$changes = @{
    2  = @("58.362.89", "  -4.66%  ")
    3  = @("2.557.04", "  -4.21%  ")
    4  = @($null, "  +0.09%  ")
    5  = @("506.17", "  -5.04%  ")
    6  = @("145.25", "  -7.35%  ")
    7  = @("0.999", "  +0.15%  ")
    8  = @("0.573", "  -3.21%  ")
    9  = @("2.568.30", "  -4.41%  ")
    10 = @("6.27", "  -5.38%  ")
    11 = @($null, "  -6.25%  ")
    12 = @("0.334", "  -5.86%  ")
    13 = @($null, "  -0.91%  ")
    14 = @("3.009.01", "  -4.05%  ")
    15 = @("58.343.36", "  -4.67%  ")
    16 = @("21.00", "  -5.36%  ")
    17 = @("0.0000136", "  -5.38%  ")
    18 = @("2.572.02", "  -3.95%  ")
    19 = @("4.52", "  -5.91%  ")
    20 = @("343.01", "  -3.83%  ")
    21 = @("10.25", "  -4.79%  ")
    22 = @("6.03", "  -4.90%  ")
    23 = @("1.00", "  +0.09%  ")
    24 = @("60.63", "  -1.71%  ")
    25 = @("0.415", "  -4.42%  ")
    26 = @("0.999", "  +0.04%  ")
    27 = @("2.685.45", "  -3.41%  ")
    28 = @($null, "  -5.85%  ")
    29 = @("0.0₃0807", "  -7.13%  ")
    30 = @("6.96", "  -6.42%  ")
    31 = @($null, "  +0.03%  ")
    32 = @("6.03", "  -2.87%  ")
    33 = @("18.69", "  -4.81%  ")
    34 = @("149.70", "  -0.34%  ")
    35 = @($null, "  -6.27%  ")
    36 = @("0.930", "  +5.11%  ")
    37 = @("3.95", "  -5.23%  ")
    38 = @($null, "  -6.66%  ")
    39 = @("0.848", "  -7.56%  ")
    40 = @("36.02", "  -2.55%  ")
    41 = @("292.12", "  -5.64%  ")
    42 = @("1.39", "  -7.86%  ")
    43 = @("3.56", "  -6.92%  ")
    44 = @("0.0992", "  -2.99%  ")
    45 = @("0.998", $null)
    46 = @("0.607", "  -6.92%  ")
    47 = @("0.0534", "  -5.96%  ")
    48 = @("19.09", "  -7.41%  ")
    49 = @($null, "  -1.08%  ")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($dVal -ne $null) {
        Set-TextValue $row 4 $dVal
    }
    if ($eVal -ne $null) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

# Rows 50 and 51 swap positions: RenderToken <-> VeChain, with refreshed data.
$ws.Cells.Item(50, 2).Value = "VeChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 50 4 "0.0227"
$ws.Cells.Item(50, 5).Value = "  -5.83%  "

$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 51 4 "4.57"
$ws.Cells.Item(51, 5).Value = "  -9.22%  "
